# Automatische test-sync: 2025-06-22 18:42:50
# Adds the new "Stageverzoek" mail-log entry as row 12 on the "Logs" sheet
# and refreshes the "Dashboard" category summary accordingly.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new inbound-mail row ---------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A12").Value = "Stageverzoek"
$logs.Range("B12").Value = "mailmind.test@zohomail.eu"
$logs.Range("C12").Value = "Ik zoek een stageplek vanaf september. Is dit mogelijk bij jullie?"
$logs.Range("D12").Value = "Sollicitatie / Vacature"
$logs.Range("F12").Value = "2025-06-22 18:42:15"
$logs.Range("G12").Value = "Nee"

# Extend the conditional-formatting ranges so row 12 is covered too.
$catRules = $logs.Range("D2:D11").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D12"))
}

$answeredRules = $logs.Range("G2:G11").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G12"))
}

# --- Dashboard sheet: "Sollicitatie / Vacature" now outnumbers "Klacht / Probleem" ---
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A3").Value = "Sollicitatie / Vacature"
$dashboard.Range("B3").Value = 2
$dashboard.Range("A4").Value = "Klacht / Probleem"
$dashboard.Range("B4").Value = 1
